# updated children task and schedules
# - re-randomised the per-trial counterbalancing data (cols A:I) for the
#   existing 24 training trials: alienID / step-count columns reshuffled and
#   the practice-length column (I) bumped from 24 to 30
# - appended 6 more training trials (rows 26:31 / trial # 25-30) using the
#   same train_dim1_2 condition label as the rest of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(1, 2, 4, 3, 4, 1, 0, 51, 30),
    @(2, 4, 3, 9, 3, 5, 0, 11, 30),
    @(3, 1, 3, 1, 3, 0, 0, 61, 30),
    @(4, 3, 2, 7, 2, 4, 0, 21, 30),
    @(5, 0, 2, 3, 2, 3, 0, 31, 30),
    @(6, 3, 4, 4, 4, 1, 0, 51, 30),
    @(7, 1, 1, 5, 1, 4, 0, 21, 30),
    @(8, 3, 3, 8, 3, 5, 0, 11, 30),
    @(9, 2, 4, 2, 4, 0, 0, 61, 30),
    @(10, 2, 3, 5, 3, 3, 0, 31, 30),
    @(11, 2, 2, 7, 2, 5, 0, 11, 30),
    @(12, 1, 0, 2, 0, 1, 0, 51, 30),
    @(13, 0, 1, 3, 1, 3, 0, 31, 30),
    @(14, 0, 3, 4, 3, 4, 0, 21, 30),
    @(15, 4, 1, 4, 1, 0, 0, 61, 30),
    @(16, 1, 3, 6, 3, 5, 0, 11, 30),
    @(17, 3, 1, 7, 1, 4, 0, 21, 30),
    @(18, 3, 3, 3, 3, 0, 0, 61, 30),
    @(19, 4, 0, 5, 0, 1, 0, 51, 30),
    @(20, 2, 1, 5, 1, 3, 0, 31, 30),
    @(21, 4, 2, 8, 2, 4, 0, 21, 30),
    @(22, 0, 2, 0, 2, 0, 0, 61, 30),
    @(23, 4, 1, 9, 1, 5, 0, 11, 30),
    @(24, 4, 4, 5, 4, 1, 0, 51, 30),
    @(25, 1, 2, 4, 2, 3, 0, 31, 30),
    @(26, 2, 0, 6, 0, 4, 0, 21, 30),
    @(27, 1, 4, 2, 4, 1, 0, 51, 30),
    @(28, 0, 0, 5, 0, 5, 0, 11, 30),
    @(29, 2, 1, 2, 1, 0, 0, 61, 30),
    @(30, 0, 4, 3, 4, 3, 0, 31, 30)
)

$numCols = 9
$arr = New-Object 'object[,]' $rows.Length, $numCols
for ($i = 0; $i -lt $rows.Length; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$lastRow = 1 + $rows.Length
$ws.Range("A2:I$lastRow").Value = $arr

# J column keeps the existing "train_dim1_2" label; fill it down into the
# newly appended rows (26:31)
for ($r = 26; $r -le $lastRow; $r++) {
    $ws.Range("J$r").Value = "train_dim1_2"
}

# selection / view state left on the newly appended block, like the source file
$ws.Range("A27:K$lastRow").Select()
